$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-apply the cell formats that the new, shorter data range needs
#        BEFORE the old trailing columns (K:P) are removed, since the
#        source cells for the copy/paste live there. ---

# D5:G5 adopt the same number-format/font as the rest of row 5 (H5),
# which leaves the old "164/fontId 6" format completely unused.
$ws.Range("H5").Copy()
$ws.Range("D5:G5").PasteSpecial(-4122) # xlPasteFormats

# E6:J6 adopt the format that used to start at M6 further along row 6;
# D6 keeps its original format.
$ws.Range("M6").Copy()
$ws.Range("E6:J6").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- 2. Update the year header row (row 4): 2007-2019 -> 2015-2021 ---
$ws.Range("D4").Value = 2015
$ws.Range("E4").Value = 2016
$ws.Range("F4").Value = 2017
$ws.Range("G4").Value = 2018
$ws.Range("H4").Value = 2019
$ws.Range("I4").Value = 2020
$ws.Range("J4").Value = 2021

# --- 3. Update row 5 (growth rate, bottom 40%) values D:J ---
$ws.Range("D5").Value = 2.2197193775563164
$ws.Range("E5").Value = 2.1235271668715399
$ws.Range("F5").Value = 2.7818537161298167
$ws.Range("G5").Value = 6.7272960584548969
$ws.Range("H5").Value = 5.1525830614767187
$ws.Range("I5").Value = 4.4774536255935971
$ws.Range("J5").Value = 4.6024666695867751

# --- 4. Update row 6 (growth rate, whole population) values D:J ---
$ws.Range("D6").Value = 2.2322863217945752
$ws.Range("E6").Value = 2.8603553109638966
$ws.Range("F6").Value = 3.113207036164539
$ws.Range("G6").Value = 6.2970593463100784
$ws.Range("H6").Value = 4.8617746111834492
$ws.Range("I6").Value = 2.6715092780025032
$ws.Range("J6").Value = 4.3694509108608912

# --- 5. Drop the now-unused trailing columns K:P (shrinks dimension to A1:J6) ---
$ws.Columns("K:P").Delete()

# --- 6. The now-numeric D:J columns get a touch wider, fixed width ---
$ws.Range("D1:J1").EntireColumn.ColumnWidth = 8.592447916666666

# --- 7. Match the saved selection/active cell ---
$ws.Range("K16").Select()
